$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(56).Insert()
$ws.Cells.Item(56, 3).Value = "update help for interface to have new positions for buttons and speed and stuff"
Write-Host "done"
